$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# --- Fix existing typos / wording ---
$ws.Range("B5").Value = 'Destruction des racines'
$ws.Range("A6").Value = 'Feuille grimpante'
$ws.Range("B8").Value = 'Placer un bouclier ou autre pour faire "Paraplui"'
$ws.Range("B12").Value = 'Cliquer pour réduire la vie du feu'

# --- Add new rows 13-19 ---
$ws.Range("A13").Value = 'Attaque distance'
$ws.Range("B13").Value = 'Cliquer pour détruire/tuer'

$ws.Range("A14").Value = 'Attaque à bout portant '
$ws.Range("B14").Value = '""'

$ws.Range("A15").Value = 'Attaque venant du ciel (pluie, bulles, cercles de fumée, animaux, objets, nourriture...)'
$ws.Range("B15").Value = 'Cliquer pour enlever'

$ws.Range("A16").Value = 'Qui foncent en ligne droite ou zigzag sur le personnage ("")'
$ws.Range("B16").Value = '""'

$ws.Range("A17").Value = 'Objets roulant : Boule de neige (montagne), tronc d''arbre (foret)'
$ws.Range("B17").Value = 'Cliquer pour le dévier les objets'

$ws.Range("A18").Value = 'Sol modifié  : Boue -> pieds qui s''enfoncent (foret), marshmallow nuage (sol collant caramel), glisse (montagne)'

$ws.Range("A19").Value = 'Vide/trou'
$ws.Range("B19").Value = 'Mettre une planche'

# --- Formatting ---
$ws.Columns.Item(1).ColumnWidth = 92.88671875

# --- Selection ---
$ws.Range("B18").Select()
